$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "ИМЯ/First Name"
$ws.Range("D1").Value = "ФАМИЛИЯ/Last Name"
$ws.Range("E1").Value = "ЛОГИН/Username"
$ws.Range("F1").Value = "НОМЕР ТЕЛ./Phone"

$ws.Rows.Item(3).Delete()
